$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 12) for year 2021, continuing the yearly time
# series that currently runs from row 2 (2011年) through row 11 (2020年).

$ws.Range("A12").Value = "2021年"
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats: reuse the same bold/bordered/centered style as the other year cells in column A

$numericValues = @{
    "B12"  = 2
    "C12"  = 161
    "D12"  = 123
    "E12"  = 4738
    "F12"  = 3
    "G12"  = 787
    "I12"  = 53
    "J12"  = 5
    "L12"  = 6400
    "M12"  = 61
    "O12"  = 44
    "P12"  = 84
    "Q12"  = 3
    "U12"  = 7
    "Y12"  = 2783
    "AA12" = 17780
    "AC12" = 119
    "AE12" = 103
    "AF12" = 234
    "AI12" = 13
    "AJ12" = 56
    "AK12" = 1315
    "AN12" = 21
    "AO12" = 16
    "AP12" = 137
    "AQ12" = 2
    "AS12" = 16
    "AT12" = 21
    "AU12" = 1
    "AV12" = 203
    "AX12" = 269
}

foreach ($addr in $numericValues.Keys) {
    $ws.Range($addr).Value = $numericValues[$addr]
}

# Columns with no reported value for 2021 still hold an explicit empty-text
# cell (matching every other year's row), rather than being left fully
# blank. Typing a bare quote enters empty text; re-applying the default
# style afterwards drops the quote-prefix formatting that Excel would
# otherwise remember for that cell.
$blankCells = @(
    "H12", "K12", "N12", "R12", "S12", "T12", "V12", "W12", "X12",
    "Z12", "AB12", "AD12", "AG12", "AH12", "AL12", "AM12", "AR12", "AW12"
)

$ws.Range("B2").Copy()
foreach ($addr in $blankCells) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats: strip the quote-prefix formatting Excel applied for the literal leading apostrophe
}
